$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.031890850384236606
$ws.Range("B1").Value = -0.031874155364268852
$ws.Range("A2").Value = -0.01048112697195098
$ws.Range("B2").Value = -0.010588592636381991
$ws.Range("A3").Value = -0.077427550191516911
$ws.Range("B3").Value = -0.077429069539474288
